$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Force text storage (matches original shared-string / text cells) instead of
# letting Excel auto-convert numeric-looking strings into numbers.
$ws.Range("B11:D12").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"

# Enterprises density (per 1000 people) - row 11
$ws.Range("B11").Value = "30.43"
$ws.Range("C11").Value = "1.37"
$ws.Range("D11").Value = "31.81"

# Employment (% of total) - row 12
$ws.Range("B12").Value = "45.68"
$ws.Range("C12").Value = "29.43"
$ws.Range("D12").Value = "75.11"

# Enterprises (% of total) - row 14
$ws.Range("B14").Value = "95.53"
$ws.Range("D14").Value = "99.83"
